$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Swap the activity note from C4 (Cedric Stephani) to B4 (Hugo Baird),
# expanding the text to describe the new section of work.
$ws.Range("C4").Value = ""
$ws.Range("B4").Value = "Gave instructions on how to download Microsoft Visual Code, ReactJS, Node.JS and MongoDB. Create starting React project base structure and push to Git. Ensure all group members can pull the project and run the react App in their browser. Assist with scheduling and delegating of stories."
